$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 5543.0557
$ws.Range("I34").Value = 1418.4
$ws.Range("J34").Value = 26166.334
$ws.Range("K34").Value = 1418.4
$ws.Range("L34").Value = 26166.334
$ws.Range("M34").Value = -1215.4
$ws.Range("N34").Value = -26572.334
$ws.Range("H36").Value = 5543.0557
$ws.Range("I36").Value = 1418.4
$ws.Range("J36").Value = 26166.334
$ws.Range("K36").Value = 1418.4
$ws.Range("L36").Value = 26166.334
$ws.Range("M36").Value = -703.4000000000001
$ws.Range("N36").Value = -27596.334
$ws.Range("H74").Value = 5039.091
$ws.Range("I74").Value = 4720
$ws.Range("J74").Value = 5422
$ws.Range("K74").Value = 4720
$ws.Range("L74").Value = 5422
$ws.Range("M74").Value = -3784
$ws.Range("N74").Value = -7294
$ws.Range("H77").Value = 5039.091
$ws.Range("I77").Value = 4720
$ws.Range("J77").Value = 5422
$ws.Range("K77").Value = 23600
$ws.Range("L77").Value = 27110
$ws.Range("M77").Value = -18920
$ws.Range("N77").Value = -36470
$ws.Range("H98").Value = 2005.2858
$ws.Range("I98").Value = 639
$ws.Range("J98").Value = 10203
$ws.Range("K98").Value = 639
$ws.Range("L98").Value = 10203
$ws.Range("M98").Value = 859
$ws.Range("N98").Value = -13199
$ws.Range("H122").Value = 2005.2858
$ws.Range("I122").Value = 639
$ws.Range("J122").Value = 10203
$ws.Range("K122").Value = 1917
$ws.Range("L122").Value = 30609
$ws.Range("M122").Value = 533
$ws.Range("N122").Value = -35509
$ws.Range("H137").Value = 2234844.5
$ws.Range("I137").Value = 2785126.5
$ws.Range("J137").Value = 1660637.4
$ws.Range("K137").Value = 8355379.5
$ws.Range("L137").Value = 4981912.199999999
$ws.Range("M137").Value = -8352829.5
$ws.Range("N137").Value = -4987012.199999999
$ws.Range("H141").Value = 358.75
$ws.Range("I141").Value = 358.75
$ws.Range("K141").Value = 1076.25
$ws.Range("M141").Value = 4103.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9966649
$ws.Range("I74").Value = 5920126.5
$ws.Range("J74").Value = 33436476
$ws.Range("K74").Value = 5920126.5
$ws.Range("L74").Value = 33436476
$ws.Range("M74").Value = -5919252.5
$ws.Range("N74").Value = -33438224
$ws.Range("H77").Value = 9966649
$ws.Range("I77").Value = 5920126.5
$ws.Range("J77").Value = 33436476
$ws.Range("K77").Value = 29600632.5
$ws.Range("L77").Value = 167182380
$ws.Range("M77").Value = -29596264.5
$ws.Range("N77").Value = -167191116
$ws.Range("H112").Value = 49000
$ws.Range("J112").Value = 49000
$ws.Range("L112").Value = 49000
$ws.Range("N112").Value = -51954
$ws.Range("H117").Value = 43996.668
$ws.Range("J117").Value = 43996.668
$ws.Range("L117").Value = 43996.668
$ws.Range("N117").Value = -53174.668
$ws.Range("H119").Value = 48500
$ws.Range("J119").Value = 48500
$ws.Range("L119").Value = 48500
$ws.Range("N119").Value = -58176
$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -48494
$ws.Range("H122").Value = 3085.1333
$ws.Range("I122").Value = 2708.3333
$ws.Range("J122").Value = 4592.3335
$ws.Range("K122").Value = 8124.999899999999
$ws.Range("L122").Value = 13777.0005
$ws.Range("M122").Value = -5674.999899999999
$ws.Range("N122").Value = -18677.0005
$ws.Range("H132").Value = 27883.46
$ws.Range("I132").Value = 40021.883
$ws.Range("K132").Value = 120065.649
$ws.Range("M132").Value = -117535.649

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1018.3611
$ws.Range("I94").Value = 920.1070999999999
$ws.Range("J94").Value = 1362.25
$ws.Range("K94").Value = 920.1070999999999
$ws.Range("L94").Value = 1362.25
$ws.Range("M94").Value = -469.1070999999999
$ws.Range("N94").Value = -2264.25
$ws.Range("H112").Value = 41000
$ws.Range("J112").Value = 41000
$ws.Range("L112").Value = 41000
$ws.Range("N112").Value = -43954
$ws.Range("H134").Value = 5287.564
$ws.Range("I134").Value = 5622.3438
$ws.Range("J134").Value = 3757.1428
$ws.Range("K134").Value = 16867.0314
$ws.Range("L134").Value = 11271.4284
$ws.Range("M134").Value = -14332.0314
$ws.Range("N134").Value = -16341.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22881
$ws.Range("J28").Value = 22881
$ws.Range("L28").Value = 22881
$ws.Range("N28").Value = -23371
$ws.Range("H107").Value = 483.7857
$ws.Range("J107").Value = 522.2
$ws.Range("L107").Value = 522.2
$ws.Range("N107").Value = -4362.2
$ws.Range("H132").Value = 13890907
$ws.Range("I132").Value = 25001250
$ws.Range("J132").Value = 2977.5
$ws.Range("K132").Value = 75003750
$ws.Range("L132").Value = 8932.5
$ws.Range("M132").Value = -75001220
$ws.Range("N132").Value = -13992.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1247.931
$ws.Range("I68").Value = 650
$ws.Range("J68").Value = 1806
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 5418
$ws.Range("M68").Value = -1139
$ws.Range("N68").Value = -7040
$ws.Range("H71").Value = 1247.931
$ws.Range("I71").Value = 650
$ws.Range("J71").Value = 1806
$ws.Range("K71").Value = 5850
$ws.Range("L71").Value = 16254
$ws.Range("M71").Value = -1794
$ws.Range("N71").Value = -24366
$ws.Range("H75").Value = 134562.94
$ws.Range("J75").Value = 84851.664
$ws.Range("L75").Value = 254554.992
$ws.Range("N75").Value = -256550.992
$ws.Range("H78").Value = 134562.94
$ws.Range("J78").Value = 84851.664
$ws.Range("L78").Value = 763664.976
$ws.Range("N78").Value = -773648.976
$ws.Range("H113").Value = 638.10345
$ws.Range("I113").Value = 599.9474
$ws.Range("J113").Value = 710.6
$ws.Range("K113").Value = 1799.8422
$ws.Range("L113").Value = 2131.8
$ws.Range("M113").Value = 370.1578
$ws.Range("N113").Value = -6471.8
$ws.Range("H131").Value = 1307.375
$ws.Range("I131").Value = 1480
$ws.Range("J131").Value = 1272.85
$ws.Range("K131").Value = 4440
$ws.Range("L131").Value = 3818.55
$ws.Range("M131").Value = 600
$ws.Range("N131").Value = -13898.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2503195
$ws.Range("I132").Value = 3574206.2
$ws.Range("J132").Value = 4169.067
$ws.Range("K132").Value = 10722618.6
$ws.Range("L132").Value = 12507.201
$ws.Range("M132").Value = -10720088.6
$ws.Range("N132").Value = -17567.201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 10113.5
$ws.Range("J104").Value = 10113.5
$ws.Range("L104").Value = 10113.5
$ws.Range("N104").Value = -17101.5
$ws.Range("H110").Value = 14198.6
$ws.Range("J110").Value = 14198.6
$ws.Range("L110").Value = 14198.6
$ws.Range("N110").Value = -22378.6
$ws.Range("H127").Value = 49950
$ws.Range("J127").Value = 49950
$ws.Range("L127").Value = 49950
$ws.Range("N127").Value = -59870
$ws.Range("H132").Value = 18532486
$ws.Range("I132").Value = 4062.5
$ws.Range("K132").Value = 12187.5
$ws.Range("M132").Value = -9657.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 223455.56
$ws.Range("I3").Value = 400360
$ws.Range("J3").Value = 2325
$ws.Range("K3").Value = 400360
$ws.Range("L3").Value = 2325
$ws.Range("M3").Value = -400246
$ws.Range("N3").Value = -2553
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676
$ws.Range("H132").Value = 5127.091
$ws.Range("I132").Value = 7666.6665
$ws.Range("J132").Value = 4174.75
$ws.Range("K132").Value = 22999.9995
$ws.Range("L132").Value = 12524.25
$ws.Range("M132").Value = -20469.9995
$ws.Range("N132").Value = -17584.25
